$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Cells.Item(7, 1).Value = 9894.6200000000008
$ws.Cells.Item(7, 2).Value = 9849.31
$ws.Cells.Item(7, 3).Value = 283.47000000000003
$ws.Cells.Item(7, 4).Value = 284.77
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = 0.46
$ws.Cells.Item(7, 7).Value = 42613.767094907409
$ws.Cells.Item(7, 8).Value = $true

# Row 8
$ws.Cells.Item(8, 1).Value = 9896.6
$ws.Cells.Item(8, 2).Value = 9894.6200000000008
$ws.Cells.Item(8, 3).Value = 282.39
$ws.Cells.Item(8, 4).Value = 282.45999999999998
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = 0.02
$ws.Cells.Item(8, 7).Value = 42614.674131944441
$ws.Cells.Item(8, 8).Value = $true

# Row 9
$ws.Cells.Item(9, 1).Value = 9930.25
$ws.Cells.Item(9, 2).Value = 9896.6
$ws.Cells.Item(9, 3).Value = 280.62
$ws.Cells.Item(9, 4).Value = 281.57
$ws.Cells.Item(9, 5).Value = $false
$ws.Cells.Item(9, 6).Value = 0.34
$ws.Cells.Item(9, 7).Value = 42615.752743055556
$ws.Cells.Item(9, 8).Value = $true

# Copy the date/time number format from G6 (existing date column style)
# onto the newly added G7:G9 cells without creating a brand-new style entry.
$ws.Range("G6").Copy()
$ws.Range("G7:G9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
